# melhorias na interface ac
# - remove os dois registros antigos de Ar Condicionado (linhas 1 e 2),
#   deixando so os aparelhos de TV no topo da lista
# - renomeia o "Ar1" restante para "ar2" e ajusta seus valores
# - adiciona um novo aparelho de A/C "ar3"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove as linhas "Ar da sala" e "Ar do quarto" - as linhas de TV sobem
# (TV da sala / TV da cozinha / TV da area passam a ser 1, 2 e 3)
$ws.Rows.Item(1).Delete()
$ws.Rows.Item(1).Delete()

# A antiga linha "Ar1" (agora linha 4) vira "ar2" com nova potencia/estado
$ws.Cells.Item(4, 1).Value = "ar2"
$ws.Cells.Item(4, 2).Value = "A/C"
$ws.Cells.Item(4, 3).Value = 17
$ws.Cells.Item(4, 4).Value = $true

# Novo aparelho de A/C "ar3"
$ws.Cells.Item(5, 1).Value = "ar3"
$ws.Cells.Item(5, 2).Value = "A/C"
$ws.Cells.Item(5, 3).Value = 30
$ws.Cells.Item(5, 4).Value = $true

# A nova linha herda a formatacao "crua" (sem estilo) da linha "ar2" acima
$ws.Range("A4:D4").Copy()
$ws.Range("A5:D5").PasteSpecial(-4122)
